$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly update: insert two new rows at the top of the "Pepino ensalada" block
# (row 338) to hold this week's prices, pushing all the previously recorded
# weeks down by two rows (the two oldest weeks at the bottom of the sheet
# simply move to the new last rows, 428-429).
$ws.Rows.Item(338).EntireRow.Insert()
$ws.Rows.Item(338).EntireRow.Insert()

# New row 338: "Primera" quality, this week's data
$ws.Cells.Item(338, 1).Value = 1
$ws.Cells.Item(338, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(338, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(338, 4).Value = 44964
$ws.Cells.Item(338, 5).Value = 15
$ws.Cells.Item(338, 6).Value = 100112043
$ws.Cells.Item(338, 7).Value = "Pepino ensalada"
$ws.Cells.Item(338, 8).Value = "Sin especificar"
$ws.Cells.Item(338, 9).Value = "Primera"
$ws.Cells.Item(338, 10).Value = 136
$ws.Cells.Item(338, 11).Value = 3500
$ws.Cells.Item(338, 12).Value = 4000
$ws.Cells.Item(338, 13).Value = 3750
$ws.Cells.Item(338, 14).Value = "$/caja 70 unidades"
$ws.Cells.Item(338, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(338, 16).Value = 54
$ws.Cells.Item(338, 17).Value = 70
$ws.Cells.Item(338, 18).Value = "Hortaliza"

# New row 339: "Segunda" quality, this week's data
$ws.Cells.Item(339, 1).Value = 1
$ws.Cells.Item(339, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(339, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(339, 4).Value = 44964
$ws.Cells.Item(339, 5).Value = 15
$ws.Cells.Item(339, 6).Value = 100112043
$ws.Cells.Item(339, 7).Value = "Pepino ensalada"
$ws.Cells.Item(339, 8).Value = "Sin especificar"
$ws.Cells.Item(339, 9).Value = "Segunda"
$ws.Cells.Item(339, 10).Value = 160
$ws.Cells.Item(339, 11).Value = 2500
$ws.Cells.Item(339, 12).Value = 3000
$ws.Cells.Item(339, 13).Value = 2750
$ws.Cells.Item(339, 14).Value = "$/caja 100 unidades"
$ws.Cells.Item(339, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(339, 16).Value = 28
$ws.Cells.Item(339, 17).Value = 100
$ws.Cells.Item(339, 18).Value = "Hortaliza"
